$d = $word.ActiveDocument

function Set-ParagraphRuns($paraIndex, [string[]]$tokens, [string]$styleName) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    # Exclude the trailing paragraph mark, then delete the paragraph's
    # content outright -- this drops any inherited run formatting (bold,
    # italics, ...) so the freshly typed text starts from a clean rPr.
    $r.MoveEnd(1, -1)
    $r.Delete()

    $p = $d.Paragraphs.Item($paraIndex)
    # Insert every token separated by a paragraph mark; each token becomes
    # its own (temporary) paragraph.
    $p.Range.InsertBefore([string]::Join("`r", $tokens))

    # Merge the inter-token paragraph marks back together. Deleting a
    # paragraph mark merges the two paragraphs it separated but leaves each
    # side's runs distinct, so we end up with N separate <w:r> runs inside a
    # single paragraph instead of Word's usual "same formatting -> merge"
    # behavior.
    for ($i = 0; $i -lt ($tokens.Length - 1); $i++) {
        $cur = $d.Paragraphs.Item($paraIndex)
        $markRange = $d.Range($cur.Range.End - 1, $cur.Range.End)
        $markRange.Delete()
    }

    $final = $d.Paragraphs.Item($paraIndex)
    $final.Style = $styleName
    return $final
}

# --- Paragraph 1: "On Pilgrimage: Russia, II" (Heading1) -> Title, split into per-token runs ---
Set-ParagraphRuns 1 @("On", " ", "Pilgrimage", ":", " ", "Russia", ",", " ", "II") "Title" | Out-Null

# --- Paragraph 2: "By Dorothy Day" (bold) -> "Dorothy Day" (Authors style, not bold), per-token runs ---
Set-ParagraphRuns 2 @("Dorothy", " ", "Day") "Authors" | Out-Null
